$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("L2").Value = 0

# Row 3
$ws.Range("L3").Value = 0

# Row 4
$ws.Range("A4").Value = 10.7
$ws.Range("H4").Value = 1.0004
$ws.Range("J4").Value = 0.98

# Row 5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0

# Row 6
$ws.Range("K6").Value = 10.6
$ws.Range("L6").Value = 0

# Row 7
$ws.Range("F7").Value = 35
$ws.Range("L7").Value = 0

# Row 8
$ws.Range("G8").Value = 147
$ws.Range("L8").Value = 0

# Row 9
$ws.Range("L9").Value = 0

# Row 10
$ws.Range("L10").Value = 0

# Row 11
$ws.Range("K11").Value = 10.9
$ws.Range("L11").Value = 0

$wb.Save()
